$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell's value to remain TEXT (matches source data which stores
# numeric-looking strings like coin prices as inline strings, not numbers), while
# resetting the cell style back to the default afterwards so no stray number-format
# style gets attached to the cell (keeps output structurally identical to source).
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) '26.658.11'
$ws.Cells.Item(2, 5).Value = '  -0.45%  '

Set-TextValue $ws.Cells.Item(3, 4) '1.596.93'
$ws.Cells.Item(3, 5).Value = '  -0.10%  '

$ws.Cells.Item(4, 5).Value = '  +0.04%  '

Set-TextValue $ws.Cells.Item(5, 4) '211.31'
$ws.Cells.Item(5, 5).Value = '  +0.32%  '

Set-TextValue $ws.Cells.Item(6, 4) '0.511'
$ws.Cells.Item(6, 5).Value = '  +1.24%  '

$ws.Cells.Item(7, 5).Value = '  +0.05%  '

$ws.Cells.Item(8, 5).Value = '  +0.25%  '

$ws.Cells.Item(9, 5).Value = '  -0.89%  '

Set-TextValue $ws.Cells.Item(10, 4) '19.66'
$ws.Cells.Item(10, 5).Value = '  +0.25%  '

Set-TextValue $ws.Cells.Item(11, 4) '0.0836'
$ws.Cells.Item(11, 5).Value = '  -0.34%  '

Set-TextValue $ws.Cells.Item(12, 4) '1.819.94'
$ws.Cells.Item(12, 5).Value = '  -0.13%  '

Set-TextValue $ws.Cells.Item(13, 4) '1.597.00'
$ws.Cells.Item(13, 5).Value = '  -0.16%  '

$ws.Cells.Item(14, 5).Value = '  -0.48%  '

$ws.Cells.Item(15, 5).Value = '  -1.42%  '

Set-TextValue $ws.Cells.Item(16, 4) '64.86'
$ws.Cells.Item(16, 5).Value = '  +2.28%  '

Set-TextValue $ws.Cells.Item(17, 4) '26.651.16'
$ws.Cells.Item(17, 5).Value = '  -0.32%  '

Set-TextValue $ws.Cells.Item(18, 4) '0.0₃0727'
$ws.Cells.Item(18, 5).Value = '  +0.02%  '

Set-TextValue $ws.Cells.Item(19, 4) '209.39'
$ws.Cells.Item(19, 5).Value = '  +0.33%  '

$ws.Cells.Item(20, 5).Value = '  -0.03%  '

$ws.Cells.Item(21, 5).Value = '  +0.38%  '

$ws.Cells.Item(22, 5).Value = '  -0.15%  '

$ws.Cells.Item(23, 5).Value = '  -1.42%  '

$ws.Cells.Item(24, 5).Value = '  +0.74%  '

Set-TextValue $ws.Cells.Item(25, 4) '146.28'
$ws.Cells.Item(25, 5).Value = '  +0.07%  '

$ws.Cells.Item(26, 5).Value = '  +0.11%  '

Set-TextValue $ws.Cells.Item(27, 4) '7.17'
$ws.Cells.Item(27, 5).Value = '  -4.18%  '

$ws.Cells.Item(28, 5).Value = '  +2.50%  '

Set-TextValue $ws.Cells.Item(29, 4) '15.28'
$ws.Cells.Item(29, 5).Value = '  +0.09%  '

$ws.Cells.Item(30, 5).Value = '  +0.88%  '

$ws.Cells.Item(31, 5).Value = '  +0.54%  '

$ws.Cells.Item(32, 5).Value = '  -0.74%  '

Set-TextValue $ws.Cells.Item(33, 4) '0.663'
$ws.Cells.Item(33, 5).Value = '  -0.52%  '

$ws.Cells.Item(34, 5).Value = '  -0.59%  '

Set-TextValue $ws.Cells.Item(35, 4) '1.298.93'
$ws.Cells.Item(35, 5).Value = '  -1.03%  '

$ws.Cells.Item(36, 5).Value = '  +0.66%  '

Set-TextValue $ws.Cells.Item(37, 4) '1.49'
$ws.Cells.Item(37, 5).Value = '  -2.09%  '

$ws.Cells.Item(38, 5).Value = '  -0.89%  '

$ws.Cells.Item(39, 5).Value = '  +2.42%  '

$ws.Cells.Item(40, 5).Value = '  -0.02%  '

$ws.Cells.Item(41, 2).Value = 'MXToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Cells.Item(41, 4) '2.20'
$ws.Cells.Item(41, 5).Value = '  +1.19%  '

$ws.Cells.Item(42, 2).Value = 'FraxShare'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Cells.Item(42, 4) '5.37'
$ws.Cells.Item(42, 5).Value = '  +2.09%  '

Set-TextValue $ws.Cells.Item(43, 4) '0.789'
$ws.Cells.Item(43, 5).Value = '  +0.15%  '

Set-TextValue $ws.Cells.Item(44, 4) '63.74'
$ws.Cells.Item(44, 5).Value = '  +1.58%  '

Set-TextValue $ws.Cells.Item(45, 4) '1.733.08'
$ws.Cells.Item(45, 5).Value = '  +0.11%  '

Set-TextValue $ws.Cells.Item(46, 4) '0.887'
$ws.Cells.Item(46, 5).Value = '  +8.54%  '

Set-TextValue $ws.Cells.Item(47, 4) '90.14'
$ws.Cells.Item(47, 5).Value = '  +1.52%  '

$ws.Cells.Item(48, 5).Value = '  +1.15%  '

Set-TextValue $ws.Cells.Item(49, 4) '0.100'
$ws.Cells.Item(49, 5).Value = '  +2.66%  '

$ws.Cells.Item(50, 5).Value = '  -0.95%  '

Set-TextValue $ws.Cells.Item(51, 4) '7.48'
$ws.Cells.Item(51, 5).Value = '  +0.90%  '
